$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above the current row 536, shifting all
# subsequent rows (old 536-610) down by one (new 537-611).
$ws.Rows.Item(536).Insert()

# Populate the newly inserted row 536 with its data. Columns that keep
# the same value as the row directly below (which held the old row-536
# data) are re-written explicitly for clarity/robustness.
$ws.Range("A536").Value = 6
$ws.Range("B536").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C536").Value = "Metropolitana"
$ws.Range("D536").Value = 45154
$ws.Range("D536").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E536").Value = 13
$ws.Range("F536").Value = 100112032
$ws.Range("G536").Value = "Zapallo italiano"
$ws.Range("H536").Value = "Sin especificar"
$ws.Range("I536").Value = "Primera"
$ws.Range("J536").Value = 760
$ws.Range("K536").Value = 14000
$ws.Range("L536").Value = 15000
$ws.Range("M536").Value = 14658
$ws.Range("N536").Value = "$/caja 50 unidades"
$ws.Range("O536").Value = "Región de Arica y Parinacota"
$ws.Range("P536").Value = 293
$ws.Range("Q536").Value = 50
$ws.Range("R536").Value = "Hortaliza"
